# AB - Game Flow
# Update Attack (column C) balance values on the Player sheet
# and move the active selection, matching the authored diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Player")

$ws.Range("C3").Value = 40
$ws.Range("C4").Value = 80
$ws.Range("C5").Value = 120
$ws.Range("C7").Value = 200
$ws.Range("C8").Value = 240
$ws.Range("C9").Value = 280
$ws.Range("C10").Value = 320
$ws.Range("C11").Value = 360
$ws.Range("C12").Value = 400

$ws.Range("G10").Select()
